$d = $word.ActiveDocument

# Replace the "prixintra" merge-field token with "prixintra_calc"
$d.Content.Find.Execute("+++=prixintra+++", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "+++=prixintra_calc+++", 2)
